$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# 1. communes -> municipalities
Replace-Text "The large distribution is part of the most important sectors at the economic level for the French territory and a fortiori for the whole world. In 2020, there are more than 44000 food sales outlets listed and at least 10900 communes in France have at least one general grocery store." "The large distribution is part of the most important sectors at the economic level for the French territory and a fortiori for the whole world. In 2020, there are more than 44000 food sales outlets listed and at least 10900 municipalities in France have at least one general grocery store."

# 2. 'Politique' paragraph -> "Politics"
Replace-Text "'Politique' is translated as 'Policy'." "Politics"

# 3. large number -> great number (political domain paragraph)
Replace-Text "The political domain engages the stability of States and the ease of exchanging with the international for all that concerns certain types of products. Indeed, large distributions work regularly with foreign countries to be able to have in their rays a large number of references and thus satisfy consumers." "The political domain engages the stability of States and the ease of exchanging with the international for all that concerns certain types of products. Indeed, large distributions work regularly with foreign countries to be able to have in their rays a great number of references and thus satisfy consumers."

# 4. large-scale retail trade paragraph rewrite
Replace-Text "The large-scale retail trade has a great influence on the economy of a country because food products are essential for all consumers." "The large distribution has necessarily a great influence on the economy of a country, since food products are part of the essentials for all consumers."

# 5. competition certainly very high -> very high; "it was obviously considered an essential business" -> "they were obviously considered essential businesses"
Replace-Text "It is a sector where competition is certainly very high, but which also attracts a large number of customers every day, regardless of the brands. The health crisis did not affect this type of business, which remained open since it was obviously considered an essential business." "It is a sector where competition is very high, but which also attracts a large number of customers every day, regardless of the brands. The health crisis did not affect this type of business, which remained open since they were obviously considered essential businesses."

# 6. brands of the sector -> brands in the sector
Replace-Text "If customers tend to desert for some of them the too large hypermarkets in favor of smaller sales spaces, it remains that the brands of the sector do not know particular difficulties." "If customers tend to desert for some of them the too large hypermarkets in favor of smaller sales spaces, it remains that the brands in the sector do not know particular difficulties."

# 7. also fight -> also to fight
Replace-Text "The clientele in the large distribution is very vast, there is no particular targeting on the part of the brands. The most important thing therefore is to be able to satisfy everyone and also fight against competitors." "The clientele in the large distribution is very vast, there is no particular targeting on the part of the brands. The most important thing therefore is to be able to satisfy everyone and also to fight against competitors."

# 8. other sectors -> the other sectors
Replace-Text "The behavior of consumers towards large-scale retail trade is different from that analyzed in other sectors. In 2020, more than 80% of customers go to the supermarket closest to their home." "The behavior of consumers towards large-scale retail trade is different from that analyzed in the other sectors. In 2020, more than 80% of customers go to the supermarket closest to their home."

# 9. Most of them are loyal to customers... -> Most loyal customers with...
Replace-Text "Most of them are loyal to customers through a loyalty card that allows them to accumulate points and earn discounts on all kinds of products." "Most loyal customers with a loyalty card that allows them to accumulate points and earn discounts on all kinds of products."

# 10. social networks sentence rewrite
Replace-Text "The large distribution wants to be in 2020 closer to the needs and expectations of consumers, with an increased presence on the net and especially on social networks. These last allow to fight more effectively against competitors by keeping a permanent interaction with Internet users." "The large distribution wants to be in 2020 closer to the needs and expectations of consumers, with an increased presence on the net and more particularly on social networks. These last ones make it possible to fight always more effectively against competitors while keeping a permanent interaction with Internet users."

# 11. ecological paragraph rewrite
Replace-Text "The large distribution, like a majority of other companies, has made these last years a big effort to adapt itself to the ecological demand of partners as well as customers. Indeed, the organic products appeared in the shelves there are some years and references do not cease increasing still today. The consumers are very in demand on this type of product which little by little, extended to all the fields of activity, food certainly, but also household products and textile." "1) The large distribution, like a majority of others companies, made these last years a big effort to adapt itself to the ecological demand of partners as well as customers. Indeed, the organic products appeared in the shelves there are some years and references do not cease increasing still today. Consumers are very in demand on this type of product which little by little, extended to all fields of activity, food certainly, but also household products and textile."

# 12. In addition, the bags -> On the other hand, bags
Replace-Text "In addition, the bags become reusable at will by the customer and many brands deliver cardboard boxes for click and collect." "On the other hand, bags become reusable at will by the customer and many brands deliver cardboard boxes for click and collect."

# 13. terms of withdrawal for technology products -> regarding technological products
Replace-Text "The terms of withdrawal for technology products or household appliances are the same as in specialized stores." "The terms of withdrawal regarding technological products or household appliances are the same as in specialized stores."

# 14. Nielsen references paragraph full rewrite
Replace-Text "According to Nielsen, sales in the large distribution sector increased by +2.5% in 2020 compared to 2019. This increase is mainly due to the health crisis and the closure of restaurants which has led consumers to buy more food products from supermarkets. In addition, online sales have also increased significantly (+20%) during this period." "According to Nielsen, sales in the large distribution sector increased by +2.5% in value and +1.7% in volume over the year 2020. This is a good performance for this sector which has been suffering from declining sales since 2014. The increase in sales is mainly due to the rise in prices (+3.8%) and not to an increase in volumes sold (-2%). The large distribution sector is thus continuing its transformation, with more and more stores being closed or converted into other types of businesses (supermarkets, hypermarkets, etc.)."
